# Auto-generated Excel COM-interop script applying the diff changes
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- sheet1 (展览) ---
$ws1.Range("F4").Value = 14
$ws1.Range("F5").Value = 500
$ws1.Range("F6").Value = 455
$ws1.Range("F8").Value = 287
$ws1.Range("F10").Value = 12984
$ws1.Range("F11").Value = 288
$ws1.Range("F12").Value = 41
$ws1.Range("F15").Value = 206
$ws1.Range("F16").Value = 163
$ws1.Range("F17").Value = 207
$ws1.Range("F18").Value = 2789
$ws1.Range("F19").Value = 53
$ws1.Range("F20").Value = 109
$ws1.Range("F21").Value = 2126
$ws1.Range("F22").Value = 162
$ws1.Range("F24").Value = 419
$ws1.Range("F26").Value = 2185
$ws1.Range("F28").Value = 1187
$ws1.Range("F29").Value = 4402
$ws1.Range("F31").Value = 4013
$ws1.Range("F32").Value = 1045
$ws1.Range("F33").Value = 2711
$ws1.Range("F34").Value = 3126
$ws1.Range("F35").Value = 109
$ws1.Range("F36").Value = 1436
$ws1.Range("F37").Value = 225
$ws1.Range("F38").Value = 804
$ws1.Range("F39").Value = 68
$ws1.Range("F40").Value = 185
$ws1.Range("F41").Value = 690
$ws1.Range("F42").Value = 1042
$ws1.Range("F43").Value = 95
$ws1.Range("F44").Value = 190
$ws1.Range("F45").Value = 431
$ws1.Range("F46").Value = 134
$ws1.Range("F47").Value = 230
$ws1.Range("F48").Value = 266

# --- sheet2 (演出) ---
$ws2.Range("F7").Value = 44
$ws2.Range("F11").Value = 15
$ws2.Range("F13").Value = 16
$ws2.Range("F18").Value = 43

# --- sheet4 (全部类型) ---
$ws4.Range("F3").Value = 501
$ws4.Range("F4").Value = 455
$ws4.Range("F6").Value = 287
$ws4.Range("F7").Value = 12984
$ws4.Range("F8").Value = 288
$ws4.Range("F11").Value = 44
$ws4.Range("F12").Value = 206
$ws4.Range("F13").Value = 163
$ws4.Range("C14").Value = "北京·Rie fu日本知名唱作歌手2024出道20周年中国巡回演唱会"
$ws4.Range("D14").Value = "奥园西路1号院4-5号楼 福浪LiveHouse"
$ws4.Range("E14").Value = "2024.05.18 20:00-05.18 22:00"
$ws4.Range("F14").Value = 45
$ws4.Range("G14").Value = 380
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=81445"
$ws4.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202401/6e9JD6401706239890264.jpeg"
$ws4.Range("C15").Value = "北京·YIYOU二次元大聚会"
$ws4.Range("D15").Value = "京开高速入口与京开高速交叉口西180米 北京双马文体创业园"
$ws4.Range("E15").Value = "2024.05.18 10:00-05.18 18:00"
$ws4.Range("F15").Value = 207
$ws4.Range("G15").Value = 55
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=83129"
$ws4.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202403/ZhTtVA3A1710812150528.png"
$ws4.Range("C16").Value = "北京·原神only3.0"
$ws4.Range("E16").Value = "2024.05.18 10:00-05.19 17:00"
$ws4.Range("F16").Value = 2789
$ws4.Range("G16").Value = 68
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=81766"
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202402/Lfxwe5PO1707120983684.jpeg"
$ws4.Range("F17").Value = 2126
$ws4.Range("F18").Value = 162
$ws4.Range("F20").Value = 419
$ws4.Range("F22").Value = 15
$ws4.Range("F23").Value = 2185
$ws4.Range("F24").Value = 1187
$ws4.Range("F27").Value = 4402
$ws4.Range("F29").Value = 4013
$ws4.Range("F30").Value = 1045
$ws4.Range("F31").Value = 2711
$ws4.Range("F32").Value = 3126
$ws4.Range("F33").Value = 109
$ws4.Range("F35").Value = 1436
$ws4.Range("F37").Value = 225
$ws4.Range("F38").Value = 804
$ws4.Range("F39").Value = 68
$ws4.Range("F40").Value = 185
$ws4.Range("F41").Value = 690
$ws4.Range("F42").Value = 43
$ws4.Range("F43").Value = 1042
$ws4.Range("F44").Value = 95
$ws4.Range("F45").Value = 190
$ws4.Range("F46").Value = 431
$ws4.Range("F47").Value = 134
$ws4.Range("F48").Value = 230
$ws4.Range("F49").Value = 266
